$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C group numbers (305020XXX -> 403110XXX, offset +98090000)
# and clear the cell style of column C (it no longer carries the
# "vertical top" alignment style - it reverts to the default/Normal style).
for ($r = 2; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 98090000
    $cell.Style = "Normal"
}

# Move the active selection to G24 (was H5)
[void]$ws.Range("G24").Select()

$wb.Save()
